$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) and "全部类型" (All types) both contain the same
# "想去人数" (number of people interested) figures in column F that were
# incremented by 1 in this update.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 127
    $ws.Range("F9").Value = 562
}
